$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the bordered/bold style from A16 down into the three new rows (A17:A19)
# so they match the existing header-style formatting used in column A for data rows.
$ws.Range("A16").Copy($ws.Range("A17:A19"))

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.91528448868924
$ws.Range("D10").Value = 1.130738473172443
$ws.Range("E10").Value = 1.183960349202367
$ws.Range("F10").Value = 0.9949748176569884
$ws.Range("G10").Value = 1.91528448868924
$ws.Range("H10").Value = 1.130738473172443
$ws.Range("I10").Value = 0.9930788932318985
$ws.Range("J10").Value = 0.6005972044409383
$ws.Range("K10").Value = 1.17902922840198
$ws.Range("L10").Value = 0.9894688374312056
$ws.Range("M10").Value = 1.91528448868924
$ws.Range("N10").Value = 1.157349411187405
$ws.Range("O10").Value = 1.30623953218026
$ws.Range("P10").Value = 1.123391536528383

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 1.406265451253938
$ws.Range("D11").Value = 0.09228280080040317
$ws.Range("E11").Value = 0.9151446579131717
$ws.Range("F11").Value = 1.156369383058451
$ws.Range("G11").Value = 1.406265451253938
$ws.Range("H11").Value = 0.09228280080040317
$ws.Range("I11").Value = 0.8229085247312959
$ws.Range("J11").Value = 1.276708738663419
$ws.Range("K11").Value = 1.054793865641898
$ws.Range("L11").Value = 0.6492028294480127
$ws.Range("M11").Value = 1.406265451253938
$ws.Range("N11").Value = 0.5037137293567874
$ws.Range("O11").Value = 0.8925155732564909
$ws.Range("P11").Value = 0.9217095314388237

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 1.38376833013707
$ws.Range("D12").Value = 0.09251911989488816
$ws.Range("E12").Value = 0.9126718756982911
$ws.Range("F12").Value = 1.159887080422128
$ws.Range("G12").Value = 1.38376833013707
$ws.Range("H12").Value = 0.09251911989488816
$ws.Range("I12").Value = 0.8229554305596779
$ws.Range("J12").Value = 1.280615927297058
$ws.Range("K12").Value = 1.055259884002806
$ws.Range("L12").Value = 0.651096048931593
$ws.Range("M12").Value = 1.38376833013707
$ws.Range("N12").Value = 0.5025954977965896
$ws.Range("O12").Value = 0.8872116015380944
$ws.Range("P12").Value = 0.9198467121179391

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 1.403056472725851
$ws.Range("D13").Value = 0.09232278934332443
$ws.Range("E13").Value = 0.9097618745324261
$ws.Range("F13").Value = 1.157611287287294
$ws.Range("G13").Value = 1.403056472725851
$ws.Range("H13").Value = 0.09232278934332443
$ws.Range("I13").Value = 0.8218383102818427
$ws.Range("J13").Value = 1.278044054420114
$ws.Range("K13").Value = 1.05507146685958
$ws.Range("L13").Value = 0.6497640817243933
$ws.Range("M13").Value = 1.403056472725851
$ws.Range("N13").Value = 0.5010423319378753
$ws.Range("O13").Value = 0.8906881059722239
$ws.Range("P13").Value = 0.9209337921468531

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.009159999999999956
$ws.Range("D14").Value = 0.06573600000000011
$ws.Range("E14").Value = 4.645183999999989
$ws.Range("F14").Value = 0.7127800000000021
$ws.Range("G14").Value = 0.009159999999999956
$ws.Range("H14").Value = 0.06573600000000011
$ws.Range("I14").Value = 1.586959999999998
$ws.Range("J14").Value = 1.310192000000003
$ws.Range("K14").Value = 0.4104959999999994
$ws.Range("L14").Value = 0.4427840000000011
$ws.Range("M14").Value = 0.009159999999999956
$ws.Range("N14").Value = 2.355459999999995
$ws.Range("O14").Value = 1.358214999999998
$ws.Range("P14").Value = 1.147911499999999

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 7.656125000000013
$ws.Range("F15").Value = 0.01
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 2.397650000000005
$ws.Range("J15").Value = 0.8912625000000018
$ws.Range("K15").Value = 0.03584999999999995
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 3.828062500000006
$ws.Range("O15").Value = 1.916531250000003
$ws.Range("P15").Value = 1.373860937500002

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.4228813996031968
$ws.Range("D16").Value = 0.428524917555197
$ws.Range("E16").Value = 4.727936168243192
$ws.Range("F16").Value = 0.4071820177408021
$ws.Range("G16").Value = 0.4228813996031968
$ws.Range("H16").Value = 0.428524917555197
$ws.Range("I16").Value = 1.810683299020801
$ws.Range("J16").Value = 0.9462085196800004
$ws.Range("K16").Value = 0.4373327578112015
$ws.Range("L16").Value = 0.4111153276928
$ws.Range("M16").Value = 0.4229070284799968
$ws.Range("N16").Value = 2.578230542899194
$ws.Range("O16").Value = 1.496631125785597
$ws.Range("P16").Value = 1.198983050918399

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 1.010127284223264
$ws.Range("D17").Value = 1.006881331918
$ws.Range("E17").Value = 0.9864321734210058
$ws.Range("F17").Value = 0.9844116600315139
$ws.Range("G17").Value = 1.010127284223264
$ws.Range("H17").Value = 1.006881331918
$ws.Range("I17").Value = 0.9897166144504399
$ws.Range("J17").Value = 0.9906846957101368
$ws.Range("K17").Value = 0.986854133907504
$ws.Range("L17").Value = 0.988719404247315
$ws.Range("M17").Value = 1.010128672423575
$ws.Range("N17").Value = 0.9966567526695027
$ws.Range("O17").Value = 0.9969631123984457
$ws.Range("P17").Value = 0.9929784122386474

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 0.6875301861625185
$ws.Range("D18").Value = 0.810245093685503
$ws.Range("E18").Value = 1.218217035534891
$ws.Range("F18").Value = 1.12651266236365
$ws.Range("G18").Value = 0.6875301861625185
$ws.Range("H18").Value = 0.810245093685503
$ws.Range("I18").Value = 0.962077573426891
$ws.Range("J18").Value = 0.9932136303389241
$ws.Range("K18").Value = 1.060436882988555
$ws.Range("L18").Value = 0.9868716132287834
$ws.Range("M18").Value = 0.6875301861625185
$ws.Range("N18").Value = 1.014231064610197
$ws.Range("O18").Value = 0.9606262444366406
$ws.Range("P18").Value = 0.9806380847162146

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.9926850836892759
$ws.Range("D19").Value = 1.254209534058356
$ws.Range("E19").Value = 1.038899482278588
$ws.Range("F19").Value = 0.9290318264355135
$ws.Range("G19").Value = 0.9926850836892759
$ws.Range("H19").Value = 1.254209534058356
$ws.Range("I19").Value = 0.9777496267929062
$ws.Range("J19").Value = 0.9673104940917976
$ws.Range("K19").Value = 0.9261600377335842
$ws.Range("L19").Value = 1.068739909352108
$ws.Range("M19").Value = 0.9927525055033046
$ws.Range("N19").Value = 1.146554508168472
$ws.Range("O19").Value = 1.053706481615433
$ws.Range("P19").Value = 1.019348249304016
